$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells AZ1:BB1 for the "3 line" GA-ACO-All chart columns ---
# (copy format from an existing header cell, A1, so the new headers get the
# same bold / centered / bordered style, then overwrite the text)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AZ1:BB1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("AZ1").Value = "Image GA-ACO All Percobaan 1"
$ws.Range("BA1").Value = "Image GA-ACO All Percobaan 2"
$ws.Range("BB1").Value = "Image GA-ACO All Percobaan 3"

# --- Add new data cells AZ2:BB6 (image paths for the new "All" chart column) ---
$ws.Range("AZ2").Value = "./tcImages/burma14_1_GA_ACO_All_10.png"
$ws.Range("BA2").Value = "./tcImages/burma14_2_GA_ACO_All_10.png"
$ws.Range("BB2").Value = "./tcImages/burma14_3_GA_ACO_All_10.png"
$ws.Range("AZ3").Value = "./tcImages/burma14_1_GA_ACO_All_50.png"
$ws.Range("BA3").Value = "./tcImages/burma14_2_GA_ACO_All_50.png"
$ws.Range("BB3").Value = "./tcImages/burma14_3_GA_ACO_All_50.png"
$ws.Range("AZ4").Value = "./tcImages/burma14_1_GA_ACO_All_100.png"
$ws.Range("BA4").Value = "./tcImages/burma14_2_GA_ACO_All_100.png"
$ws.Range("BB4").Value = "./tcImages/burma14_3_GA_ACO_All_100.png"
$ws.Range("AZ5").Value = "./tcImages/lin318_1_GA_ACO_All_10.png"
$ws.Range("BA5").Value = "./tcImages/lin318_2_GA_ACO_All_10.png"
$ws.Range("BB5").Value = "./tcImages/lin318_3_GA_ACO_All_10.png"
$ws.Range("AZ6").Value = "./tcImages/lin318_1_GA_ACO_All_50.png"
$ws.Range("BA6").Value = "./tcImages/lin318_2_GA_ACO_All_50.png"
$ws.Range("BB6").Value = "./tcImages/lin318_3_GA_ACO_All_50.png"

# --- Update existing numeric cell values (new simulation run results) ---
$ws.Range("J2").Value = 49.70749645258855
$ws.Range("K2").Value = 51.22634287474276
$ws.Range("L2").Value = 41.74882913758471
$ws.Range("M2").Value = 47.56088948830534
$ws.Range("N2").Value = 5.090359295161472
$ws.Range("Q2").Value = 31.88252949105588
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 31.45623383762054
$ws.Range("V2").Value = 31.88252949105587
$ws.Range("W2").Value = 31.74043093991076
$ws.Range("X2").Value = 0.2461219102652595
$ws.Range("Y2").Value = 0.0003
$ws.Range("AB2").Value = 0.0003
$ws.Range("AI2").Value = 0.0005
$ws.Range("AJ2").Value = 0.0004666666666666667
$ws.Range("J3").Value = 37.06866231056024
$ws.Range("K3").Value = 36.23408339635454
$ws.Range("L3").Value = 39.89511019355628
$ws.Range("M3").Value = 37.73261863349035
$ws.Range("N3").Value = 1.918699507876747
$ws.Range("P3").Value = 31.88252949105588
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 31.22691510942754
$ws.Range("V3").Value = 31.88252949105588
$ws.Range("W3").Value = 31.66399136384643
$ws.Range("X3").Value = 0.3785191397177077
$ws.Range("AD3").Value = 0.0009
$ws.Range("AF3").Value = 0.0009000000000000001
$ws.Range("J4").Value = 36.61066473569531
$ws.Range("K4").Value = 31.80719333670128
$ws.Range("L4").Value = 43.38553809791235
$ws.Range("M4").Value = 37.26779872343631
$ws.Range("N4").Value = 5.817077071927895
$ws.Range("T4").Value = 32.11184821924887
$ws.Range("V4").Value = 31.22691510942754
$ws.Range("W4").Value = 31.74043093991077
$ws.Range("X4").Value = 0.4592609834208278
$ws.Range("Y4").Value = 0.0023
$ws.Range("AB4").Value = 0.0023
$ws.Range("AE4").Value = 0.0019
$ws.Range("AF4").Value = 0.001766666666666667
$ws.Range("AG4").Value = 0.004
$ws.Range("AJ4").Value = 0.004
$ws.Range("J5").Value = 563066.2639079965
$ws.Range("K5").Value = 574205.9909832703
$ws.Range("L5").Value = 572895.2970148403
$ws.Range("M5").Value = 570055.850635369
$ws.Range("N5").Value = 6088.531976947872
$ws.Range("Q5").Value = 49168.33743272909
$ws.Range("R5").Value = 49199.85415701829
$ws.Range("S5").Value = 27.29428387851803
$ws.Range("V5").Value = 49215.61251916289
$ws.Range("W5").Value = 49218.04237919758
$ws.Range("X5").Value = 75.5136558155243
$ws.Range("Y5").Value = 0.0063
$ws.Range("AA5").Value = 0.0065
$ws.Range("AB5").Value = 0.006433333333333333
$ws.Range("AC5").Value = 0.0645
$ws.Range("AD5").Value = 0.06469999999999999
$ws.Range("AE5").Value = 0.0636
$ws.Range("AF5").Value = 0.06426666666666665
$ws.Range("AG5").Value = 0.07140000000000001
$ws.Range("AH5").Value = 0.0718
$ws.Range("AI5").Value = 0.0713
$ws.Range("AJ5").Value = 0.07149999999999999
$ws.Range("J6").Value = 534594.0100370105
$ws.Range("K6").Value = 533135.3539705164
$ws.Range("L6").Value = 545755.2663368067
$ws.Range("M6").Value = 537828.2101147779
$ws.Range("N6").Value = 6903.664580531461
$ws.Range("O6").Value = 48452.08503658375
$ws.Range("P6").Value = 48786.55204087113
$ws.Range("R6").Value = 48818.08319887259
$ws.Range("S6").Value = 382.7390946938857
$ws.Range("U6").Value = 48452.08503658375
$ws.Range("V6").Value = 48272.74602537625
$ws.Range("W6").Value = 48646.81452704097
$ws.Range("X6").Value = 500.6885130664281
$ws.Range("Y6").Value = 0.0284
$ws.Range("Z6").Value = 0.0289
$ws.Range("AA6").Value = 0.029
$ws.Range("AB6").Value = 0.02876666666666667
$ws.Range("AC6").Value = 0.1577
$ws.Range("AD6").Value = 0.1604
$ws.Range("AE6").Value = 0.1606
$ws.Range("AF6").Value = 0.1595666666666667
$ws.Range("AG6").Value = 0.1859
$ws.Range("AH6").Value = 0.1896
$ws.Range("AI6").Value = 0.1912
$ws.Range("AJ6").Value = 0.1889
$ws.Range("Q6").Value = 49215.61251916289
